# Tumorkonferenzen-IG observations-summary sheet:
#  - Row 23 ("TNMYSmbol") had a typo in both its Profile (A) and Name (B)
#    cells; fix it to "TNMYSymbol".
#  - Row 24 was an exact duplicate of row 23 (same Code/VS/Types, already
#    spelled "TNMYSymbol") — once row 23's typo is fixed the two rows are
#    fully identical, so remove row 24 entirely. Excel shifts rows 25-29
#    up into 24-28 and shrinks the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "TNMYSymbol"
$ws.Range("B23").Value = "TNMYSymbol"

$ws.Rows.Item(24).Delete()
